$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
# Refresh the account-statement data table (rows 16-22): the "Periodo Mora"
# (E) column is re-sorted from descending (1909..1903) to ascending
# (1903..1909), "Valor Mora" (F) stays tied to its period (33125 for every
# period except 1909, which keeps 25396), and "Salario Basico" (G) is
# updated from 0 to 828116 for every row.

$periodos = @("1903", "1904", "1905", "1906", "1907", "1908", "1909")
$valorMora = @(33125, 33125, 33125, 33125, 33125, 33125, 25396)
$salarioBasico = @(828116, 828116, 828116, 828116, 828116, 828116, 828116)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valorMora[$i]
    $ws.Range("G$row").Value = $salarioBasico[$i]
}
